$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the next day's gold-price row (row 42), matching the pattern of
# the preceding rows: column A holds the date as text, column B reuses
# the same "price of gold" text as the previous day's row (price unchanged
# since the latest data point).
$lastRow = 41
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "26-10-2025"
$ws.Cells.Item($newRow, 2).Value = $ws.Cells.Item($lastRow, 2).Text

# Match formatting of the row being appended after.
$ws.Range("A" + $lastRow + ":B" + $lastRow).Copy() | Out-Null
$ws.Range("A" + $newRow + ":B" + $newRow).PasteSpecial(-4122) | Out-Null
